$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 74 (revised quarterly figures) -------------------
$ws.Range("B74").Value = 1510
$ws.Range("C74").Value = 1061
$ws.Range("D74").Value = 449
$ws.Range("E74").Value = 44
$ws.Range("F74").Value = 319
$ws.Range("G74").Value = 210
$ws.Range("H74").Value = 666
$ws.Range("I74").Value = 508
$ws.Range("J74").Value = 3257

# --- Append new row 75 for period 01-04-2021 -------------------------------
# A75 holds the text "01-04-2021". Assigning that literal string straight to
# .Value would make Excel's input-parser treat it as a date and stamp the
# cell with an explicit (date) number format, which the source workbook does
# not have (column A cells are plain shared-string text with the default
# style). Build the text via a formula and paste its *value* in instead, so
# it lands as plain text with no number-format override, matching the
# original column A cells exactly.
$ws.Range("Z1").Formula = "=""01-04-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A75").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("B75").Value = 1654
$ws.Range("C75").Value = 1173
$ws.Range("D75").Value = 481
$ws.Range("E75").Value = 15
$ws.Range("F75").Value = 334
$ws.Range("G75").Value = 222
$ws.Range("H75").Value = 700
$ws.Range("I75").Value = 463
$ws.Range("J75").Value = 3389
